# Weekly update: insert a new week's worth of "Perejil" (parsley) price
# records for Vega Monumental Concepcion, right before the old row 66.
# This pushes every following row down by two and adds a brand-new pair
# of rows (Primera / Segunda) dated 44460, while the two rows that used
# to fall off the end of the range reappear at the bottom (102/103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at position 66 (shifts 66:101 down to 68:103).
$ws.Rows.Item(66).Insert()
$ws.Rows.Item(66).Insert()

# The row that used to be 66 is now at 68 (and 67 -> 69). Clone their full
# contents into the two new blank rows so every column matches the
# existing "Primera"/"Segunda" template for this market/product.
$ws.Range("A68:R68").Copy()
$ws.Range("A66:R66").PasteSpecial()

$ws.Range("A69:R69").Copy()
$ws.Range("A67:R67").PasteSpecial()

# Stamp the new week's date onto the freshly inserted rows.
$ws.Range("D66").Value2 = 44460
$ws.Range("D67").Value2 = 44460
